# Revised schedule, vcs video link
#
# The only real content changes in this commit are:
#  1. Row 59 (M59): the "Task 4 demo" video link is swapped for the newer
#     panopto recording (the "vcs video link" mentioned in the commit msg).
#  2. Rows 60-65 and 69: the old "Get repo" / "TBA" / GitHub Classroom
#     invite-link cells (columns F and J) are removed — those assignments
#     no longer show a "get repo" step.
#
# Removing the last references to "Get repo", "TBA", the old panopto link,
# and the two classroom.github.com links causes those shared strings to
# disappear and every later shared-string index to shift down — which is
# why the diff shows so many <v> index changes even though the *text*
# in those cells is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New panopto video link for the Task 4 demo / version-control session.
$ws.Range("M59").Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=7cb26d66-a26e-4b24-89b1-acaf01321690"

# 2. Drop the obsolete "Get repo" / "TBA" / GitHub Classroom cells.
$ws.Range("F60").ClearContents()
$ws.Range("J60").ClearContents()

$ws.Range("F61").ClearContents()
$ws.Range("J61").ClearContents()

$ws.Range("F62").ClearContents()
$ws.Range("J62").ClearContents()

$ws.Range("F63").ClearContents()
$ws.Range("J63").ClearContents()

$ws.Range("F64").ClearContents()
$ws.Range("J64").ClearContents()

$ws.Range("F65").ClearContents()
$ws.Range("J65").ClearContents()

$ws.Range("F69").ClearContents()
$ws.Range("J69").ClearContents()

# Keep the active-cell selection in sync with the author's saved state.
$ws.Range("K59").Select()
